$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight B7, D7, E7 with red font color (new style using red font)
$ws.Range("B7").Font.Color = 255
$ws.Range("D7").Font.Color = 255
$ws.Range("E7").Font.Color = 255

# Add average row (row 10): sum of rows 2-9 divided by 8
$ws.Range("B10").Formula = "=SUM(B2:B9)/8"
$ws.Range("C10").Formula = "=SUM(C2:C9)/8"
$ws.Range("D10").Formula = "=SUM(D2:D9)/8"
$ws.Range("E10").Formula = "=SUM(E2:E9)/8"

# Update view: zoom and active selection
$ws.Application.ActiveWindow.Zoom = 265
$ws.Range("C10").Select()
